# "hot topic list Jan 2019.xlsx" update -- add two new hot-topic rows and
# widen column A so the longer topic names fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: new hot topic "THz frequencies" owned by "Angel?" with description.
$ws.Range("A3").Value = "THz frequencies: communication meets positioning"
$ws.Range("B3").Value = "Angel?"
$ws.Range("C3").Value = "The next frontier after mmWave seems to be THz frequencies, where signals can simultaneously serve to communicate (short range) and to position with high resolution"

# Row 4: new hot topic "blockchain for comms", owner still to be filled in ("???").
$ws.Range("A4").Value = "blockchain for comms"
$ws.Range("B4").Value = "???"

# Move the cursor to the next empty row, like Excel leaves it after data entry.
[void]$ws.Range("A5").Select()

# Column A needs to be a lot wider to show the new, longer topic names; B/C stay
# essentially the same width (just nudged to their natural character widths).
$ws.Columns.Item(1).ColumnWidth = 39.5
$ws.Columns.Item(2).ColumnWidth = 20.333333333333336
$ws.Columns.Item(3).ColumnWidth = 48
